# "Made adjustments to date for exercise tables"
#
# The single original "Sheet1" (exercise log/table) is renamed to
# "Exercise Table" and three new tracking sheets are inserted in front of
# it: "Squats", "Standing_Lunges" and "Dumbbell_Curls" - each a small
# Date/Reps log for one exercise. A few values on the Exercise Table are
# also corrected.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet -------------------------------------------
$exerciseTable = $wb.Worksheets.Item(1)
$exerciseTable.Name = "Exercise Table"

# --- Fix up a handful of values on the Exercise Table ---------------------
# Row 3 = Squats: sets/reps corrected, a 3rd-set column added
$exerciseTable.Range("B3").Value = 12
$exerciseTable.Range("C3").Value = 76
$exerciseTable.Range("D3").Value = 7

# Row 4 = Reverse Leg Lifts: the "10" moves from the 2-set column to the
# 3-set column
$exerciseTable.Range("C4").ClearContents()
$exerciseTable.Range("B4").Value = 10

# Row 6 = Dumbbell Curls: values corrected
$exerciseTable.Range("C6").ClearContents()
$exerciseTable.Range("B6").Value = 12
$exerciseTable.Range("D6").Value = 10

# Row 10 = Sit Ups: drop the misc value
$exerciseTable.Range("D10").ClearContents()

# --- Add the three new per-exercise date/reps tracking sheets ------------
# Worksheets.Add() always inserts the new sheet at the far left, so add
# them in reverse of the desired final left-to-right order:
# Squats, Standing_Lunges, Dumbbell_Curls, Exercise Table

$dumbbellCurls = $wb.Worksheets.Add()
$dumbbellCurls.Name = "Dumbbell_Curls"
$dumbbellCurls.Range("A1").Value = "Date"
$dumbbellCurls.Range("A1").Font.Bold = $true
$dumbbellCurls.Range("A1").Font.Color = 3937500
$dumbbellCurls.Range("B1").Value = "Dumbbell_Curls"
$dumbbellCurls.Range("B1").Font.Bold = $true
$dumbbellCurls.Range("B1").Font.Color = 3937500
$dumbbellCurls.Range("A2").Value = 43528.272893518515
$dumbbellCurls.Range("A2").Font.Bold = $true
$dumbbellCurls.Range("A2").NumberFormat = "m/d/yyyy h:mm"
$dumbbellCurls.Range("B2").Value = 12
$dumbbellCurls.Range("C7").Value = 6

$standingLunges = $wb.Worksheets.Add()
$standingLunges.Name = "Standing_Lunges"
$standingLunges.Range("A1").Value = "Date"
$standingLunges.Range("A1").Font.Bold = $true
$standingLunges.Range("A1").Font.Color = 3937500
$standingLunges.Range("B1").Value = "Standing_Lunges"
$standingLunges.Range("B1").Font.Bold = $true
$standingLunges.Range("B1").Font.Color = 3937500
$standingLunges.Range("A2").Value = 43528.925520833334
$standingLunges.Range("A2").Font.Bold = $true
$standingLunges.Range("A2").NumberFormat = "m/d/yyyy h:mm"
$standingLunges.Range("B2").Value = 6
$standingLunges.Range("D3").Value = 10
$standingLunges.Range("D6").Value = 10

$squats = $wb.Worksheets.Add()
$squats.Name = "Squats"
$squats.Range("A1").Value = "Date"
$squats.Range("A1").Font.Bold = $true
$squats.Range("A1").Font.Color = 3937500
$squats.Range("B1").Value = "Squats"
$squats.Range("B1").Font.Bold = $true
$squats.Range("B1").Font.Color = 3937500
$squats.Range("A2").Value = 43558
$squats.Range("A2").Font.Bold = $true
$squats.Range("A2").NumberFormat = "m/d/yyyy"
$squats.Range("B2").Value = 10
